$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---

# Row 2
$ws.Range("C2").Value = "not fixed"
$ws.Range("D2").Value = "Open"
$ws.Range("F2").Value = "2023-08-03 14:48:21"

# Row 7
$ws.Range("D7").Value = "Open"
$ws.Range("F7").Value = "2023-08-03 16:16:08"

# Row 8
$ws.Range("D8").Value = "Closed"
$ws.Range("F8").Value = "2023-08-03 13:03:56"

# --- Append new rows 11-15, matching formatting of column A (style s="1") ---
$ws.Range("A10").Copy()
$ws.Range("A11:A15").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "good bye"
$ws.Range("C11").Value = "good bye"
$ws.Range("D11").Value = "Closed"
$ws.Range("E11").Value = "2023-08-03 13:06:27"
$ws.Range("F11").Value = "2023-08-03 13:06:50"

# Row 12
$ws.Range("A12").Value = 21
$ws.Range("B12").Value = "auto refresh"
$ws.Range("C12").Value = "auto refresh"
$ws.Range("D12").Value = "In Progress"
$ws.Range("E12").Value = "2023-08-03 14:18:12"
$ws.Range("F12").Value = "2023-08-03 14:56:56"

# Row 13
$ws.Range("A13").Value = 22
$ws.Range("B13").Value = "dddd"
$ws.Range("C13").Value = "dddd"
$ws.Range("D13").Value = "Closed"
$ws.Range("E13").Value = "2023-08-03 14:57:06"
$ws.Range("F13").Value = "2023-08-03 16:06:03"

# Row 14
$ws.Range("A14").Value = 100
$ws.Range("B14").Value = "godjflkdaj"
$ws.Range("C14").Value = "godjflkdaj"
$ws.Range("D14").Value = "Closed"
$ws.Range("E14").Value = "2023-08-03 15:58:29"
$ws.Range("F14").Value = "2023-08-03 16:05:52"

# Row 15
$ws.Range("A15").Value = 321
$ws.Range("B15").Value = "weqweq"
$ws.Range("C15").Value = "weqweq"
$ws.Range("D15").Value = "Closed"
$ws.Range("E15").Value = "2023-08-03 16:15:50"
$ws.Range("F15").Value = "2023-08-03 16:16:20"

Write-Host "Edit complete"
